$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.225.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.796.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.90%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4516'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3708'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.01'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("E10").Value = '  +3.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07570'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.64%  '
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.289'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.487'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.794.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.17%  '
$ws.Range("E17").Value = '  +3.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06746'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.369'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.215.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.413'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.363'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.999.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.239'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09409'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.798'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2369'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.61%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02342'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06304'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.206'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6570'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.368'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.481'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.205'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6082'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.822'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("E48").Value = '  +3.74%  '
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("E50").Value = '  +2.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.159'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.57%  '
